# Dynamic environment url implementation / code cleanup:
# the login-URL cell (C1) used to carry a live hyperlink to a hard-coded
# environment URL. That URL is now supplied dynamically elsewhere, so the
# hyperlink (and the text it displayed) is removed from the sheet,
# leaving the cell blank but keeping its existing formatting/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("C1")
$cell.Hyperlinks.Delete()
$cell.ClearContents()
